# "adding averages and more checks"
#
# The "LAST UPDATE" recheck date on the Training Dashboard sheet moves
# forward 8 days, from 08-Sep-2025 to 16-Sep-2025, which pulls every
# "PERIOD TO EXPIRE" day-count (column H) down by 8 for each training row.
#
# The header row's bold font (sitting on the dark navy fill) is also
# recolored to white so it is legible, and the title row is normalized to
# use that same bold/white look.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$oldDate = "08-Sep-2025"
$newDate = "16-Sep-2025"
$deltaDays = 8

for ($row = 3; $row -le 27; $row++) {
    $updateCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE
    $periodCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE

    if ($updateCell.Value2 -eq $oldDate) {
        # Force the cell to stay plain text (it already is) instead of
        # letting Excel auto-coerce the "dd-mmm-yyyy"-looking string into
        # a date serial when it is written back.
        $updateCell.NumberFormat = "@"
        $updateCell.Value = $newDate

        $periodCell.Value = $periodCell.Value2 - $deltaDays
    }
}

# Header row (row 2) and title row (row 1): bold, white font text.
$titleRange = $ws.Range("A1")
$headerRange = $ws.Range("A2:K2")

$titleRange.Font.Bold = $true
$titleRange.Font.Size = 11
$titleRange.Font.Color = 16777215

$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
